# Generate Report for Handback
# Adds a new handback row (e4be68ea-63ed-4675-832a-2bdfae7519d3) alongside the
# existing row, whose identifying GUID is renamed from
# 49cc0501-330e-4578-af8d-66c66e9a4998 to de3a9987-ddc8-4b72-b143-47b529f267d8,
# to the Overview / zh-cn / de-de worksheets (and their tables).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$guid1 = "de3a9987-ddc8-4b72-b143-47b529f267d8"
$guid2 = "e4be68ea-63ed-4675-832a-2bdfae7519d3"
$hash1 = "e9d9be27629cd6194fe1f88de8377b5885e33fb0"
$hash2 = "0580dd24c88766be5f5a1ad14bcb96ba9116dec3"

$oltestBase = "https://github.com/OpenLocalizationTestOrg/oltest/blob/6f641072cdfbcf3672e1f72b239f4e29a32378e7/e2e/"
$zhcnBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/2fe6b1cb847d81bb9365909cd2301f41fd9a740c/e2e/"
$dedeBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/72e119c50a7b9246cf3fa924c02f83381ada73d0/e2e/"

# ---------------------------------------------------------------------------
# 1) Rename the GUID that belongs to the already-present row (the "handoff"
#    that this commit regenerated) everywhere it shows up, and refresh its
#    "Latest HO Xliff Generate Date" / datetime stamps.
# ---------------------------------------------------------------------------

$wsOverview.Range("A2").Value = "$guid1.md"
$wsOverview.Range("G2").Value = "2016-08-12 09:15:13"
$wsOverview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("A2").Value = "$guid1.md"
$wsZh.Range("I2").Value = "$guid1.md"
$wsZh.Range("G2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("J2").Value = "$guid1.$hash1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-12 09:14:59"
$wsZh.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("K2").Value = "2016-08-12 09:15:30"
$wsZh.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("A2").Value = "$guid1.md"
$wsDe.Range("I2").Value = "$guid1.md"
$wsDe.Range("G2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("J2").Value = "$guid1.$hash1.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-12 09:15:13"
$wsDe.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("K2").Value = "2016-08-12 09:15:39"
$wsDe.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Hyperlinks that target the renamed GUID need refreshing too. Rebuild each
# sheet's hyperlink collection from scratch so relationship ids stay tight
# (rId2, rId3, ...) exactly like a freshly generated report would.
$wsOverview.Hyperlinks.Delete()
$wsZh.Hyperlinks.Delete()
$wsDe.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2) Grow each table by one row for the newly handed-back file.
# ---------------------------------------------------------------------------

$wsOverview.ListObjects.Item(1).ListRows.Add() | Out-Null
$wsZh.ListObjects.Item(1).ListRows.Add() | Out-Null
$wsDe.ListObjects.Item(1).ListRows.Add() | Out-Null

# ---------------------------------------------------------------------------
# 3) Overview sheet - row 3
# ---------------------------------------------------------------------------

$wsOverview.Range("A3").Value = "$guid2.md"
$wsOverview.Range("C3").Value = "'.md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-12 09:15:13"
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# 4) zh-cn sheet - row 3
# ---------------------------------------------------------------------------

$wsZh.Range("B3").Value = "'.md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-12 09:14:59"
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("J3").Value = "$guid2.$hash2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-12 09:15:30"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L3").Value = "'"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = "'"
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = "'"

# ---------------------------------------------------------------------------
# 5) de-de sheet - row 3
# ---------------------------------------------------------------------------

$wsDe.Range("B3").Value = "'.md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-12 09:15:13"
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("J3").Value = "$guid2.$hash2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-12 09:15:39"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L3").Value = "'"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = "'"
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = "'"

# ---------------------------------------------------------------------------
# 6) Hyperlinks (re-created fresh -> rId2, rId3, rId4, rId5 ...)
# ---------------------------------------------------------------------------

$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "$oltestBase$guid1.md", "", "", "e2e\$guid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "$oltestBase$guid2.md", "", "", "e2e\$guid2.md") | Out-Null

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$oltestBase$guid1.md", "", "", "$guid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$zhcnBase$guid1.md", "", "", "$guid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$oltestBase$guid2.md", "", "", "$guid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$zhcnBase$guid2.md", "", "", "$guid2.md") | Out-Null

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$oltestBase$guid1.md", "", "", "$guid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$dedeBase$guid1.md", "", "", "$guid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$oltestBase$guid2.md", "", "", "$guid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$dedeBase$guid2.md", "", "", "$guid2.md") | Out-Null

"Report regenerated for $guid2"
